$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.361.85'
$ws.Range("D3").Value = '2.454.64'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.43%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = '2.453.40'
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.93%  '
$ws.Range("E11").Value = '  -1.97%  '
$ws.Range("E12").Value = '  -6.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.69%  '
$ws.Range("D14").Value = '2.904.83'
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '68.262.57'
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.72%  '
$ws.Range("D18").Value = '2.469.29'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.68'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.74%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.79%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.09%  '
$ws.Range("D28").Value = '2.579.28'
$ws.Range("E28").Value = '  -1.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.19%  '
$ws.Range("D30").Value = '0.0₃0826'
$ws.Range("E30").Value = '  -6.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +124.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '429.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.88%  '
$ws.Range("E35").Value = '  -3.59%  '
$ws.Range("E36").Value = '  -4.20%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '156.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.27%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.109'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.304'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.02%  '
$ws.Range("E44").Value = '  -5.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").Value = '  -5.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '134.45'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.99%  '
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.479'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.562'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.67%  '
